$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$cr = [char]13
$newText = $cr + "建置ｔｅｓｔ：這樣可以同時測試同一類型的檔案是否一直都能正常運行"
$tr.InsertAfter($newText) | Out-Null
